# The "Sentiment" column (A) on sheet "sentHisLarge" used to hold the text
# labels "Positive" / "Negative". The commit replaces those labels with the
# equivalent numeric codes 1 (Positive) and 0 (Negative) for rows 2-31
# (rows 2-16 were "Positive", rows 17-31 were "Negative").
#
# Re-writing these cells as numbers automatically drops "Positive"/"Negative"
# from the shared-string table (they become unused) and shifts every other
# shared-string index down, which is exactly what the target workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sentHisLarge")

# Rows 2-16: "Positive" -> 1
$ws.Range("A2:A16").Value = 1

# Rows 17-31: "Negative" -> 0
$ws.Range("A17:A31").Value = 0

# Match the author's final cursor position/selection on the sheet.
$ws.Range("A30").Select()
